$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Julio de 2020 a las 16:27"

# Malaui and Maldivas swapped rank (row 109/110) due to updated case counts,
# so their country names need to be swapped between the two rows.
$ws.Range("A109").Value = "Malaui"
$ws.Range("A110").Value = "Maldivas"

# Update numeric statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) for the rows whose
# figures changed in this data refresh.

# Row 4
$ws.Range("B4").Value = 3964355
$ws.Range("C4").Value = 2926
$ws.Range("D4").Value = 1851125
$ws.Range("E4").Value = 1969315
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 81
$ws.Range("H4").Value = 143915

# Row 6
$ws.Range("B6").Value = 1170636
$ws.Range("C6").Value = 15719
$ws.Range("D6").Value = 735061
$ws.Range("E6").Value = 407246
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 230
$ws.Range("H6").Value = 28329

# Row 21
$ws.Range("B21").Value = 203557
$ws.Range("C21").Value = 70
$ws.Range("D21").Value = 187800
$ws.Range("E21").Value = 6582
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 9175

# Row 26
$ws.Range("B26").Value = 97159
$ws.Range("C26").Value = 2466
$ws.Range("D26").Value = 64950
$ws.Range("E26").Value = 28259
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 81
$ws.Range("H26").Value = 3950

# Row 30
$ws.Range("B30").Value = 78166
$ws.Range("C30").Value = 20
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = 5646

# Row 62
$ws.Range("B62").Value = 21605
$ws.Range("C62").Value = 352
$ws.Range("D62").Value = 14047
$ws.Range("E62").Value = 7067
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 9
$ws.Range("H62").Value = 491

# Row 63
$ws.Range("B63").Value = 21442
$ws.Range("C63").Value = 327
$ws.Range("D63").Value = 14599
$ws.Range("E63").Value = 6136
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 12
$ws.Range("H63").Value = 707

# Row 70
$ws.Range("B70").Value = 14152
$ws.Range("C70").Value = 54
$ws.Range("D70").Value = 8836
$ws.Range("E70").Value = 4957
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 359

# Row 80
$ws.Range("B80").Value = 9412
$ws.Range("C80").Value = 163
$ws.Range("D80").Value = 4940
$ws.Range("E80").Value = 4040
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 432

# Row 82
$ws.Range("B82").Value = 9049
$ws.Range("C82").Value = 15
$ws.Range("D82").Value = 8138
$ws.Range("E82").Value = 656
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 255

# Row 91
$ws.Range("B91").Value = 6967
$ws.Range("C91").Value = 46
$ws.Range("D91").Value = 5683
$ws.Range("E91").Value = 1227
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 57

# Row 101
$ws.Range("B101").Value = 4290
$ws.Range("C101").Value = 119
$ws.Range("D101").Value = 2397
$ws.Range("E101").Value = 1776
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 4
$ws.Range("H101").Value = 117

# Row 109
$ws.Range("B109").Value = 3045
$ws.Range("C109").Value = 53
$ws.Range("D109").Value = 1180
$ws.Range("E109").Value = 1801
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = 64

# Row 110
$ws.Range("B110").Value = 2999
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 2369
$ws.Range("E110").Value = 615
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 15

# Row 137
$ws.Range("B137").Value = 1366
$ws.Range("C137").Value = 22
$ws.Range("D137").Value = 57
$ws.Range("E137").Value = 1302
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 3
$ws.Range("H137").Value = 7
